# Updated cryptos list on Fri Sep 27 22:49:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Pct($val) {
    return "  " + $val + "  "
}

# Set a cell's value while forcing it to remain a text value (Excel would
# otherwise silently re-interpret numeric-looking strings, e.g. "7.70", as
# numbers and drop the trailing zero). ClearFormats() afterwards removes the
# temporary "@" text-format style so no stray style is left behind.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "65.947.04"
$ws.Range("E2").Value = Pct("+1.26%")

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.707.58"
$ws.Range("E3").Value = Pct("+2.63%")

# Row 4 - TetherUSD
$ws.Range("E4").Value = Pct("-0.05%")

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "609.38"
$ws.Range("E5").Value = Pct("+2.00%")

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "158.24"
$ws.Range("E6").Value = Pct("+1.28%")

# Row 7 - USDC
$ws.Range("E7").Value = Pct("-0.04%")

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.589"

# Row 9 - Dogecoin
$ws.Range("E9").Value = Pct("+5.39%")

# Row 10 - Toncoin
$ws.Range("E10").Value = Pct("+3.61%")

# Row 11 - Cardano
$ws.Range("E11").Value = Pct("+0.25%")

# Row 12 - TRON
$ws.Range("E12").Value = Pct("+1.08%")

# Row 13 - Avalanche
Set-TextValue $ws.Range("D13") "30.54"
$ws.Range("E13").Value = Pct("+4.48%")

# Row 14 - ShibaInu
$ws.Range("E14").Value = Pct("+8.78%")

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.192.65"
$ws.Range("E15").Value = Pct("+2.49%")

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "65.846.61"
$ws.Range("E16").Value = Pct("+1.21%")

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.710.24"
$ws.Range("E17").Value = Pct("+1.69%")

# Row 18 - Chainlink
Set-TextValue $ws.Range("D18") "12.74"
$ws.Range("E18").Value = Pct("+1.20%")

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "4.91"
$ws.Range("E19").Value = Pct("+1.96%")

# Rows 20 and 21 swap positions: BitcoinCash <-> Uniswap (with new values)
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D20") "7.70"
$ws.Range("E20").Value = Pct("+4.78%")

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D21") "360.18"
$ws.Range("E21").Value = Pct("+1.93%")

# Row 22 - Dai
$ws.Range("E22").Value = Pct("-0.10%")

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "70.95"
$ws.Range("E23").Value = Pct("+3.71%")

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D24") "9.90"
$ws.Range("E24").Value = Pct("+3.59%")

# Row 25 - PEPE
$ws.Range("E25").Value = Pct("+12.59%")

# Row 26 - SuiNetwork
$ws.Range("E26").Value = Pct("-1.56%")

# Row 27 - Fetch.AI
Set-TextValue $ws.Range("D27") "1.70"
$ws.Range("E27").Value = Pct("+3.49%")

# Row 28 - Kaspa
$ws.Range("E28").Value = Pct("+3.84%")

# Row 29 - Aptos
$ws.Range("E29").Value = Pct("+4.10%")

# Row 30 - PancakeSwap
$ws.Range("E30").Value = Pct("+5.12%")

# Row 31 - Bittensor
Set-TextValue $ws.Range("D31") "545.90"
$ws.Range("E31").Value = Pct("+6.87%")

# Row 32 - Binance-PegBSC-USD
$ws.Range("E32").Value = Pct("+0.05%")

# Row 33 - ImmutableX
$ws.Range("E33").Value = Pct("+2.46%")

# Row 34 - RenderToken
Set-TextValue $ws.Range("D34") "6.79"
$ws.Range("E34").Value = Pct("+6.84%")

# Row 35 - NEARProtocol
Set-TextValue $ws.Range("D35") "5.45"
$ws.Range("E35").Value = Pct("-2.81%")

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = Pct("+2.16%")

# Row 37 - EthereumClassic
Set-TextValue $ws.Range("D37") "20.95"
$ws.Range("E37").Value = Pct("+3.03%")

# Row 38 - Monero
Set-TextValue $ws.Range("D38") "163.39"
$ws.Range("E38").Value = Pct("-0.05%")

# Row 39 - Stacks
$ws.Range("E39").Value = Pct("-0.32%")

# Row 41 - Aave
Set-TextValue $ws.Range("D41") "173.45"
$ws.Range("E41").Value = Pct("+4.64%")

# Row 42 - USDe
$ws.Range("E42").Value = Pct("+0.01%")

# Row 43 - OKB
Set-TextValue $ws.Range("D43") "42.58"
$ws.Range("E43").Value = Pct("+0.64%")

# Row 44 - Filecoin
Set-TextValue $ws.Range("D44") "4.21"
$ws.Range("E44").Value = Pct("+2.82%")

# Row 45 - Hedera
$ws.Range("E45").Value = Pct("+0.32%")

# Row 46 - InjectiveProtocol
Set-TextValue $ws.Range("D46") "23.62"
$ws.Range("E46").Value = Pct("+2.07%")

# Row 47 - dogwifhat
$ws.Range("E47").Value = Pct("+4.18%")

# Row 48 - VeChain
$ws.Range("E48").Value = Pct("+4.45%")

# Row 49 - Mantle
$ws.Range("E49").Value = Pct("+1.36%")

# Row 50 - EnergySwap
Set-TextValue $ws.Range("D50") "21.13"
$ws.Range("E50").Value = Pct("+8.35%")

# Row 51 - Stellar
Set-TextValue $ws.Range("D51") "0.0992"
$ws.Range("E51").Value = Pct("+0.88%")
